$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New monthly row: 01-10-2021. Force text so the date-like label isn't
# auto-converted to a date serial (matches the rest of column A, which is
# stored as plain text), then drop the temporary formatting again so the
# cell keeps the workbook's default (unstyled) look.
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "01-10-2021"
$ws.Range("A35").ClearFormats()

$ws.Range("B35").Value = 6
$ws.Range("C35").Value = 5.1
$ws.Range("D35").Value = 6.9
$ws.Range("E35").Value = 4.9
$ws.Range("F35").Value = 4.8
